$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column D (Category), so Name's old column D
# (Category) through F shift right to E through G, making room for "Last Name"
# in the new D column.
$ws.Range("D1").EntireColumn.Insert()

# Update header row: split "Name" into "First Name" / "Last Name"
$ws.Range("C1").Value = "First Name"
$ws.Range("D1").Value = "Last Name"

# Split full names in column C into First/Last Name (C/D) for each data row
$lastRow = $ws.UsedRange.Rows.Count
for ($r = 2; $r -le $lastRow; $r++) {
    $fullName = $ws.Cells.Item($r, 3).Value2
    $parts = $fullName -split ' ', 2
    $first = $parts[0]
    $last = if ($parts.Length -gt 1) { $parts[1] } else { "" }
    $ws.Cells.Item($r, 3).Value = $first
    $ws.Cells.Item($r, 4).Value = $last
}
